# Update the "Contenu du stage" breakdown (rows 16-20) on the Worksheet
# sheet: student counts move from an all-C# 30/0/0/0 split to a mixed
# 6/20/0/2/2 split (C#, COBOL, C++, ASSEMBLEUR, ANDROID), and the matching
# percentage labels in column G are refreshed to reflect the new split.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# The G column holds literal percentage text (e.g. "20 %"), not a numeric
# percentage. Pre-format as text so Excel doesn't silently convert the
# "66.67 %" / "6.67 %" strings we assign below into percentage numbers.
$ws.Range("G16:G20").NumberFormat = "@"

# Row 16: C#
$ws.Range("E16").Value = 6
$ws.Range("G16").Value = "20 %"

# Row 17: COBOL
$ws.Range("E17").Value = 20
$ws.Range("G17").Value = "66.67 %"

# Row 18: C++ (unchanged values, left as-is)

# Row 19: ASSEMBLEUR
$ws.Range("E19").Value = 2
$ws.Range("G19").Value = "6.67 %"

# Row 20: ANDROID
$ws.Range("E20").Value = 2
$ws.Range("G20").Value = "6.67 %"

# Rows 21-23 (JEE, DELPHI, PHP5) and the "Type entreprise" block
# (rows 25-28) are unchanged.

Write-Output "Contenu du stage breakdown updated"
